$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.740.81'
$ws.Range("D3").Value = '1.725.03'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.31'
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9985'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").Value = '  -1.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2589'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06179'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '1.722.96'
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '15.81'
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6025'
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.83'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9985'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '26.558.05'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9980'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007135'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").Value = '1.944.58'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.409'
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.509'
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.053'
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.75'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.21'
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +2.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '106.32'
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  -2.33%  '
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07920'
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.664'
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04526'
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.593'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9976'
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6169'
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9280'
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("E38").Value = '  +2.80%  '
$ws.Range("E39").Value = '  +1.20%  '
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.605'
$ws.Range("E42").Value = '  +3.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.83'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3824'
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.769'
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05356'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.909'
$ws.Range("E48").Value = '  +2.87%  '
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.245'
$ws.Range("E50").Value = '  +2.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.45'
$ws.Range("E51").Value = '  +0.83%  '
